# Informatii_necesare_setare_cont_demo_1.xlsx — data refresh
# ---------------------------------------------------------------
# 1) A handful of placeholder names / e-mail addresses were
#    corrected on three sheets ("Cont adminstrator", "Receptie",
#    "Angajati").
# 2) The conditional-format rule used on "Domenii" / "Servicii"
#    (cells A2:A3 / A5) had its font reset from a blue, underlined
#    style to plain black text.
# 3) The workbook now opens with the "Receptie" tab active/selected
#    (cell B4) instead of "Angajati".

$wb = $excel.ActiveWorkbook

# --- 1a. "Cont adminstrator" — business name + admin e-mail ------
$wsAdmin = $wb.Worksheets.Item("Cont adminstrator")
$wsAdmin.Range("B15").Value = "Scoala particulara Mario M"
$wsAdmin.Range("C15").Value = "marioscaola1@automation.33mail.com"

# --- 1b. "Receptie" — receptionist e-mails ------------------------
$wsReceptie = $wb.Worksheets.Item("Receptie")
$wsReceptie.Range("B2").Value = "steluta2za1@staffcalendis.33mail.com"
$wsReceptie.Range("B3").Value = "steluta3131za@staffcalendis.33mail.com"
$wsReceptie.Range("B4").Value = "steluta441za@staffcalendis.33mail.com"

# --- 1c. "Angajati" — employee e-mails ----------------------------
$wsAngajati = $wb.Worksheets.Item("Angajati")
$wsAngajati.Range("B2").Value = "elenaz114a1@staffcalendis.33mail.com"
$wsAngajati.Range("B3").Value = "komornicza11124@staffcalendis.33mail.com"
$wsAngajati.Range("B4").Value = "ovidiusz11a34@staffcalendis.33mail.com"
$wsAngajati.Range("B5").Value = "sdrosea1sz4@staffcalendis.33mail.com"

# --- 2. Conditional-format font: blue/underlined -> plain black --
$wsDomenii = $wb.Worksheets.Item("Domenii")
$cond = $wsDomenii.Range("A2:A3").FormatConditions.Item(1)
$cond.Font.Underline = $false
$cond.Font.Color = 0

# --- 3. Move the active tab/selection to "Receptie"!B4 -----------
$wsReceptie.Activate() | Out-Null
$wsReceptie.Range("B4").Select() | Out-Null
